# Generate Report for Handback
# Update timestamp values in the three worksheets to reflect the latest
# handoff/handback xliff generation/report times for the
# 62f0fa30-cf3b-4c68-901f-9e9469c72187.md item.

$wb = $excel.ActiveWorkbook

# "Overview" sheet - Latest HO Xliff Generate Date for the third row (row 4)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-09-01 04:49:35"

# "zh-cn" sheet - Correspond Handoff Datetime / Correspond Handback DateTime for row 4
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-09-01 04:49:31"
$wsZhCn.Range("K4").Value = "2016-09-01 04:49:49"

# "de-de" sheet - Correspond Handback DateTime for row 4
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K4").Value = "2016-09-01 04:49:57"
